$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 7.741029
$ws.Range("H2").Value = 23.223087
$ws.Range("I2").Value = 0.4930486933812723
$ws.Range("J2").Value = 0.4930486933812723
$ws.Range("Q2").Value = 2.275343877057
$ws.Range("R2").Value = 20.478094893513
$ws.Range("S2").Value = 0.4930486933812723
$ws.Range("T2").Value = 0.4930486933812723

# Row 3
$ws.Range("I3").Value = 0.0194007766416684
$ws.Range("J3").Value = 0.0194007766416684
$ws.Range("S3").Value = 0.0194007766416684
$ws.Range("T3").Value = 0.0194007766416684

# Row 4
$ws.Range("G4").Value = 7.654706000000001
$ws.Range("I4").Value = 0.4875505299770593
$ws.Range("J4").Value = 0.4875505299770593
$ws.Range("S4").Value = 0.4875505299770593
$ws.Range("T4").Value = 0.4875505299770593
